$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "Team Members" paragraph
#   " Schmid Tosha, Bonsignori Evan, Brown Darryle, Schnibben David, Norris Tyler"
#   -> " Schmid Tosha, Bonsignori Evan, Brown Darryle" + _GoBack bookmark + ", Norris Tyler"
#   (i.e. remove ", Schnibben David" and leave a _GoBack bookmark at the split point,
#    which is where Word leaves it after the last edit made to the document)
# ------------------------------------------------------------------

$rng1 = $d.Content
$rng1.Find.Execute("Brown Darryle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng1.End

# Insert the _GoBack bookmark at the split point BEFORE deleting the trailing
# text, so the run is cleanly divided in two (runs with identical formatting
# that are not separated by a bookmark get coalesced back together on save).
$bmRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$rng2 = $d.Content
$rng2.Find.Execute(", Schnibben David", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Delete()

# ------------------------------------------------------------------
# Edit 2: "- Partially impleme" + _GoBack bookmark + "nted XTEA as mentioned in Milestones."
#   -> single run "- Partially implemented XTEA as mentioned in Milestones."
#   (the stray _GoBack bookmark that was splitting this run is removed, since
#    the bookmark now belongs at the Team Members edit made above)
# ------------------------------------------------------------------

$rng3 = $d.Content
$rng3.Find.Execute("- Partially impleme" + "nted XTEA as mentioned in Milestones.", $true, $false, $false, $false, $false, $true, 1, $false, "- Partially implemented XTEA as mentioned in Milestones.", 2) | Out-Null
